$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Convert the SmartScore values on row 4 (Julieta Hernandez) from text to
#    real numbers, keeping the same numeric value.
# ---------------------------------------------------------------------------
$ws.Range("I4").Value = 0.612
$ws.Range("L4").Value = 0.573
$ws.Range("O4").Value = 0.496
$ws.Range("R4").Value = 0.643
$ws.Range("U4").Value = 0.626
$ws.Range("X4").Value = 0.5679999999999999
$ws.Range("AA4").Value = 0.695
$ws.Range("AD4").Value = 0.6860000000000001
$ws.Range("AG4").Value = 0.645

# ---------------------------------------------------------------------------
# 2) Append a new row (row 5) with the Streamlit submission for
#    "Fernanda Adamaris".
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "Fernanda Adamaris_20251113_204405"
$ws.Range("B5").Value = ""
$ws.Range("C5").Value = "Fernanda Adamaris"
$ws.Range("D5").Value = 20
$ws.Range("E5").Value = "Female"
$ws.Range("F5").Value = "2025-11-13 20:44:06"
$ws.Range("G5").Value = "{`n  ""portion"": 0.4,`n  ""diet"": 0.14285714285714285,`n  ""salt"": 0.2,`n  ""fat"": 0.4,`n  ""natural"": 0.2,`n  ""convenience"": 0.2,`n  ""price"": 0.4`n}"

$ws.Range("H5").Value = "Nongshim Neoguri Spicy Seafood"
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "0.600"
$ws.Range("J5").Value = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"

$ws.Range("K5").Value = "Maruchan Ramen Sabor Pollo"
$ws.Range("L5").NumberFormat = "@"
$ws.Range("L5").Value = "0.563"
$ws.Range("M5").Value = "Sabor clásico, económico, alto en sodio, no saludable, nostálgico"

$ws.Range("N5").Value = "Nissin Chow Mein Teriyaki Beef"
$ws.Range("O5").NumberFormat = "@"
$ws.Range("O5").Value = "0.539"
$ws.Range("P5").Value = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"

$ws.Range("Q5").Value = "Kraft Macaroni & Cheese Dinner"
$ws.Range("R5").NumberFormat = "@"
$ws.Range("R5").Value = "0.643"
$ws.Range("S5").Value = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"

$ws.Range("T5").Value = "Velveeta Original Shells & Cheese (microwave cups)"
$ws.Range("U5").NumberFormat = "@"
$ws.Range("U5").Value = "0.582"
$ws.Range("V5").Value = "Muy cremoso, porción individual, rápido, salado, ideal para niños"

$ws.Range("W5").Value = "Annie’s Shells & White Cheddar"
$ws.Range("X5").NumberFormat = "@"
$ws.Range("X5").Value = "0.569"
$ws.Range("Y5").Value = "Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños"

$ws.Range("Z5").Value = "Wild Planet Wild Tuna Pasta Salad"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "0.624"
$ws.Range("AB5").Value = "Sabor fresco, buena proteína, saludable, porción algo pequeña"

$ws.Range("AC5").Value = "StarKist Chicken Creations (Chicken Salad)"
$ws.Range("AD5").NumberFormat = "@"
$ws.Range("AD5").Value = "0.611"
$ws.Range("AE5").Value = "Portátil, saludable, fácil, buena textura, sabor suave"

$ws.Range("AF5").Value = "Jack Link’s Beef Jerky Original"
$ws.Range("AG5").NumberFormat = "@"
$ws.Range("AG5").Value = "0.611"
$ws.Range("AH5").Value = "Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña"
